$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (many values look numeric, e.g. "10.90" or "0.000009918", and would
# otherwise be normalized/reformatted by Excel's automatic type detection).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.488.82"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.864.62"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "324.30"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "0.4556"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Value = "0.3830"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").Value = "0.07821"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "0.9880"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "21.54"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "1.846.05"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "6.899"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "5.625"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "0.06928"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "86.57"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "0.000009918"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "16.64"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "28.494.75"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "5.245"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "10.90"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("D24").Value = "2.096"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "2.086.30"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "153.65"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "19.10"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "5.669"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "1.927"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").Value = "117.29"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "0.09275"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "0.9052"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").Value = "5.252"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "1.318"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "3.292"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").Value = "0.05699"
$ws.Range("D37").Value = "1.137"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "0.02050"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "7.663"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "0.5549"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "0.1767"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "9.640"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").Value = "0.07089"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "11.52"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "0.5233"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "1.135"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "2.108"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "1.811"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").Value = "111.69"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").Value = "2.421"
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  -0.22%  "
